$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the BOURL / CCURL / CWURL values from test18 -> test4uk environment
$ws.Range("A2").Value = "https://test4uk.cliotest.com/backoffice/control/main"
$ws.Range("C2").Value = "https://test4uk.cliotest.com/cabicentral/control/main"
$ws.Range("D2").Value = "https://test4uk.cliotest.com/warehouse/control/main"

# Update vhostTarget / clockServerTarget values
$ws.Range("F2").Value = "virtual_cabitest4uk"
$ws.Range("G2").Value = "test4uk"
